$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) cells.
# NumberFormat is forced to text before assignment so that numeric-looking
# strings (e.g. "322.48") are not auto-coerced into Number cells by Excel,
# then the style is reset to Normal so no visible formatting/style is altered.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.300.55"
Set-TextValue "E2" "  +1.22%  "
Set-TextValue "D3" "1.902.97"
Set-TextValue "E3" "  +1.06%  "
Set-TextValue "E4" "  +0.11%  "
Set-TextValue "D5" "322.48"
Set-TextValue "E5" "  -2.33%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  +0.05%  "
Set-TextValue "D7" "0.4720"
Set-TextValue "E7" "  +2.86%  "
Set-TextValue "D8" "0.4031"
Set-TextValue "E8" "  -1.27%  "
Set-TextValue "D9" "0.08025"
Set-TextValue "E9" "  +0.71%  "
Set-TextValue "D10" "0.9927"
Set-TextValue "E10" "  +0.03%  "
Set-TextValue "D11" "22.63"
Set-TextValue "E11" "  +4.73%  "
Set-TextValue "D12" "1.896.61"
Set-TextValue "E12" "  +1.62%  "
Set-TextValue "D13" "5.858"
Set-TextValue "E13" "  -0.80%  "
Set-TextValue "D14" "7.053"
Set-TextValue "E14" "  -0.18%  "
Set-TextValue "D15" "89.26"
Set-TextValue "E15" "  +1.23%  "
Set-TextValue "E16" "  +0.14%  "
Set-TextValue "D17" "0.06622"
Set-TextValue "E17" "  +0.70%  "
Set-TextValue "D19" "17.56"
Set-TextValue "E19" "  +0.97%  "
Set-TextValue "D20" "1.001"
Set-TextValue "E20" "  +0.04%  "
Set-TextValue "D21" "29.317.95"
Set-TextValue "E21" "  +1.39%  "
Set-TextValue "D22" "5.512"
Set-TextValue "E22" "  +1.89%  "
Set-TextValue "D23" "11.46"
Set-TextValue "E23" "  -0.32%  "
Set-TextValue "D24" "2.203"
Set-TextValue "E24" "  +0.17%  "
Set-TextValue "D25" "2.113.73"
Set-TextValue "E25" "  +1.05%  "
Set-TextValue "D26" "154.30"
Set-TextValue "E26" "  -1.55%  "
Set-TextValue "D27" "19.69"
Set-TextValue "E27" "  +0.87%  "
Set-TextValue "D28" "6.035"
Set-TextValue "E28" "  +10.31%  "
Set-TextValue "D29" "2.088"
Set-TextValue "E29" "  +0.50%  "
Set-TextValue "D30" "117.50"
Set-TextValue "E30" "  +0.03%  "
Set-TextValue "D31" "1.067"
Set-TextValue "E31" "  +4.74%  "
Set-TextValue "D32" "0.09488"
Set-TextValue "E32" "  +1.69%  "
Set-TextValue "D33" "1.410"
Set-TextValue "E33" "  +0.69%  "
Set-TextValue "E34" "  +1.09%  "
Set-TextValue "D35" "5.359"
Set-TextValue "E35" "  +1.52%  "
Set-TextValue "D36" "0.06065"
Set-TextValue "E36" "  +0.47%  "
Set-TextValue "D37" "0.02244"
Set-TextValue "E37" "  +0.79%  "
Set-TextValue "D38" "1.175"
Set-TextValue "E38" "  -0.26%  "
Set-TextValue "E39" "  -3.07%  "
Set-TextValue "D40" "0.5803"
Set-TextValue "E40" "  +0.46%  "
Set-TextValue "D41" "2.502"
Set-TextValue "E41" "  +12.11%  "
Set-TextValue "D42" "0.1828"
Set-TextValue "E42" "  +0.24%  "
Set-TextValue "D43" "10.06"
Set-TextValue "E43" "  -0.06%  "
Set-TextValue "D44" "0.07799"
Set-TextValue "E44" "  +4.05%  "
Set-TextValue "D45" "1.275"
Set-TextValue "E45" "  +1.27%  "
Set-TextValue "D46" "12.14"
Set-TextValue "E46" "  +1.43%  "
Set-TextValue "D47" "0.5487"
Set-TextValue "E47" "  +0.72%  "
Set-TextValue "E48" "  +0.23%  "
Set-TextValue "D49" "113.46"
Set-TextValue "E49" "  +2.05%  "
Set-TextValue "D50" "43.61"
Set-TextValue "E50" "  -3.44%  "
Set-TextValue "D51" "0.2891"
Set-TextValue "E51" "  +2.91%  "
